$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The edit reorders the data rows (2-17) of the weekly Fruta/Hortaliza
# price sheet. Only columns D (Fecha), M (Volumen), N (Precio minimo),
# O (Precio maximo), P (Precio promedio ponderado), R (Origen) and
# S (Precio $/Kg) vary row to row - every other column is constant
# across the block, so only those seven columns need to be rewritten.

# New values per row (2..17), taken from the target layout.
$rows = @(
    @{ D = 44232; M = 200; N = 3000; O = 3000; P = 3000; R = "Provincia de Curicó"; S = 1500 },
    @{ D = 44978; M = 500; N = 3000; O = 3000; P = 3000; R = "Provincia de Curicó"; S = 1500 },
    @{ D = 44231; M = 150; N = 3400; O = 3400; P = 3400; R = "Provincia de Curicó"; S = 1700 },
    @{ D = 44188; M = 150; N = 3000; O = 3400; P = 3240; R = "Provincia de Linares"; S = 1620 },
    @{ D = 44237; M = 100; N = 3600; O = 4000; P = 3800; R = "Provincia de Curicó"; S = 1900 },
    @{ D = 44582; M = 380; N = 5000; O = 5000; P = 5000; R = "Provincia de Curicó"; S = 2500 },
    @{ D = 44980; M = 250; N = 4000; O = 4000; P = 4000; R = "Provincia de Curicó"; S = 2000 },
    @{ D = 44194; M = 120; N = 3000; O = 3000; P = 3000; R = "Provincia de Linares"; S = 1500 },
    @{ D = 44174; M = 200; N = 3200; O = 3200; P = 3200; R = "Provincia de Curicó"; S = 1600 },
    @{ D = 44208; M = 85;  N = 3000; O = 3000; P = 3000; R = "Provincia de Linares"; S = 1500 },
    @{ D = 44236; M = 300; N = 3600; O = 4000; P = 3800; R = "Provincia de Curicó"; S = 1900 },
    @{ D = 44238; M = 300; N = 3600; O = 4000; P = 3800; R = "Provincia de Curicó"; S = 1900 },
    @{ D = 44586; M = 250; N = 5000; O = 5000; P = 5000; R = "Provincia de Curicó"; S = 2500 },
    @{ D = 44168; M = 170; N = 8000; O = 8000; P = 8000; R = "Provincia de Linares"; S = 4000 },
    @{ D = 44617; M = 90;  N = 6500; O = 6500; P = 6500; R = "Provincia de Curicó"; S = 3250 },
    @{ D = 44533; M = 150; N = 4000; O = 4000; P = 4000; R = "Provincia de Curicó"; S = 2000 }
)

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = 2 + $i
    $data = $rows[$i]
    $ws.Cells.Item($r, 4).Value = $data.D    # D: Fecha
    $ws.Cells.Item($r, 13).Value = $data.M   # M: Volumen
    $ws.Cells.Item($r, 14).Value = $data.N   # N: Precio minimo
    $ws.Cells.Item($r, 15).Value = $data.O   # O: Precio maximo
    $ws.Cells.Item($r, 16).Value = $data.P   # P: Precio promedio ponderado
    $ws.Cells.Item($r, 18).Value = $data.R   # R: Origen
    $ws.Cells.Item($r, 19).Value = $data.S   # S: Precio $/Kg
}
